$d = $word.ActiveDocument

function Find-ParagraphStartingWith($doc, $prefix) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($prefix)) {
            return $p
        }
    }
    return $null
}

# 1) Remove the "Continuum optimisation-based refinement..." item and the
#    "PRIME3D local refinement mode..." item right after it (they are
#    consecutive list paragraphs, so delete them together as one range).
$pContinuum = Find-ParagraphStartingWith $d "Continuum optimisation-based refinement"
$pPrime3D   = $pContinuum.Next()
$d.Range($pContinuum.Range.Start, $pPrime3D.Range.End).Delete()

# 2) Remove the "Preparing method paper for Protein Science: ..." item.
$pPaper = Find-ParagraphStartingWith $d "Preparing method paper for Protein Science"
$pPaper.Range.Delete()

# 3) Remove the "Think about and develop a down-scaling strategy..." item.
$pDownscale = Find-ParagraphStartingWith $d "Think about and develop a"
$pDownscale.Range.Delete()

# 4) In the final todo item, drop the old "Approach for coarse orientation
#    search..." sentence but keep the _GoBack bookmark where it sits, then
#    append the new sentence after the (now first) bookmark.
$pLast = Find-ParagraphStartingWith $d "Approach for coarse orientation search"
$oldText = "Approach for coarse orientation search in initial stage (say 200 projection directions). These can be searched at low-res (say 20-30 A) and the most promising hits will be subjected to fine-grained search. All improving solution candidates will be assigned weights after search is done. "
$d.Range($pLast.Range.Start, $pLast.Range.Start + $oldText.Length).Delete()

$insertPoint = $d.Range($pLast.Range.End - 1, $pLast.Range.End - 1)
$insertPoint.InsertAfter("Method for diversifying the continuous refinement (based on search history???)")
